# Finished the static lists of movie titles
#
# 1) Clean up several paragraphs that had their text artificially split
#    across multiple runs (wrapped with w:proofErr spellStart/spellEnd
#    bookmarks from Word's spell checker) by replacing each paragraph's
#    content with a single clean run (preserving paragraph formatting
#    such as list numbering where present).
# 2) Append a new "New Releases" section with its own heading and movie
#    list at the end of the document.

$d = $word.ActiveDocument

function Set-ParaXml($doc, $index, $innerXml) {
    $p = $doc.Paragraphs($index)
    $rng = $p.Range
    $full = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($full)
}

$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# --- 1) Merge split runs / drop stray proofErr spell-check bookmarks ---

Set-ParaXml $d 6 '<w:r><w:t>Bodies Bodies Bodies</w:t></w:r>'
Set-ParaXml $d 7 '<w:r><w:t>Brahmastra Part One: Shiva</w:t></w:r>'
Set-ParaXml $d 10 '<w:r><w:t>Cuando Sea Joven</w:t></w:r>'
Set-ParaXml $d 12 '<w:r><w:t>Dio: Dreamers Never Die</w:t></w:r>'
Set-ParaXml $d 14 '<w:r><w:t>Dongalunnaru Jagratta</w:t></w:r>'
Set-ParaXml $d 21 '<w:r><w:t>Moonage Daydream</w:t></w:r>'
$nomadlandXml = $listPPr + '<w:r><w:t>Nomadland</w:t></w:r>'
Set-ParaXml $d 37 $nomadlandXml
$benHurXml = $listPPr + '<w:r><w:t>Ben-Hur</w:t></w:r>'
Set-ParaXml $d 98 $benHurXml

# --- 2) Append the "New Releases" section at the end of the document ---

$newReleasesXml = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>New Releases</w:t>
  </w:r>
</w:p>
<w:p><w:r><w:t>Bullet Proof</w:t></w:r></w:p>
<w:p><w:r><w:t>The Munsters</w:t></w:r></w:p>
<w:p><w:r><w:t>Thor: Love and Thunder</w:t></w:r></w:p>
<w:p><w:r><w:t>Sound of Metal</w:t></w:r></w:p>
<w:p><w:r><w:t>Lost Illusions</w:t></w:r></w:p>
<w:p><w:r><w:t>Black Phone</w:t></w:r></w:p>
<w:p><w:r><w:t>Vengeance</w:t></w:r></w:p>
<w:p><w:r><w:t>The Reef: Stalked</w:t></w:r></w:p>
<w:p><w:r><w:t>Wrong Place</w:t></w:r></w:p>
<w:p><w:r><w:t>The Forgiven</w:t></w:r></w:p>
<w:p><w:r><w:t>Elvis</w:t></w:r></w:p>
<w:p><w:r><w:t>Lightyear</w:t></w:r></w:p>
<w:p><w:r><w:t>Where the Crawdads Sing</w:t></w:r></w:p>
<w:p><w:r><w:t>Happening</w:t></w:r></w:p>
<w:p>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Mrs. Harris Goes to Paris</w:t>
  </w:r>
</w:p>
<w:p><w:r><w:t>Paradise Highway</w:t></w:r></w:p>
<w:p><w:r><w:t>Minions: The Rise of Gru</w:t></w:r></w:p>
<w:p><w:r><w:t>American Carnage</w:t></w:r></w:p>
<w:p><w:r><w:t>The Phantom of the Open</w:t></w:r></w:p>
<w:p><w:r><w:t>Murder at Yellowstone City</w:t></w:r></w:p>
<w:p><w:r><w:t>Mr. Malcom\u2019s List</w:t></w:r></w:p>
<w:p><w:r><w:t>Frank &amp; Penelope</w:t></w:r></w:p>
<w:p><w:r><w:t>Jurassic World Dominion</w:t></w:r></w:p>
<w:p><w:r><w:t>Vivo</w:t></w:r></w:p>
<w:p><w:r><w:t>Nitram</w:t></w:r></w:p>
<w:p><w:r><w:t>Crimes of the Future</w:t></w:r></w:p>
<w:p><w:r><w:t>Hot Seat</w:t></w:r></w:p>
<w:p><w:r><w:t>Firestarter</w:t></w:r></w:p>
<w:p><w:r><w:t>Sonic the Hedgehog 2</w:t></w:r></w:p>
<w:p><w:r><w:t>Last Seen Alive</w:t></w:r></w:p>
<w:p>
  <w:r>
    <w:t>Harry Potter 20</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>th</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Anniversary</w:t>
  </w:r>
</w:p>
<w:p><w:r><w:t>Downton Abbey: A New Era</w:t></w:r></w:p>
<w:p><w:r><w:t>White Elephant</w:t></w:r></w:p>
<w:p><w:r><w:t>Charm City Kings</w:t></w:r></w:p>
<w:p><w:r><w:t>Cow</w:t></w:r></w:p>
'@

$newReleasesXml = $newReleasesXml.Replace([char]0x5c + 'u2019', [char]0x2019)

$endRng = $d.Content
$endRng.Collapse(0)
$fullPkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newReleasesXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRng.InsertXML($fullPkg)
